$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9, pushing the existing data rows 9-19 down to 10-20.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new observation.
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(9, 3).Value = "Ñuble"
$ws.Cells.Item(9, 4).Value = 44775
$ws.Cells.Item(9, 4).NumberFormat = $ws.Cells.Item(10, 4).NumberFormat
$ws.Cells.Item(9, 5).Value = 16
$ws.Cells.Item(9, 6).Value = 100112037
$ws.Cells.Item(9, 7).Value = "Cebollín"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 100
$ws.Cells.Item(9, 11).Value = 8000
$ws.Cells.Item(9, 12).Value = 8000
$ws.Cells.Item(9, 13).Value = 8000
$ws.Cells.Item(9, 14).Value = "`$/docena de atados"
$ws.Cells.Item(9, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(9, 16).Value = 2667
$ws.Cells.Item(9, 17).Value = 3
$ws.Cells.Item(9, 18).Value = "Hortaliza"
